$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 22.69183541765118
$ws.Range("C2").Value = 11.10752750488409
$ws.Range("E2").Value = 10.33249327558904
$ws.Range("F2").Value = 46.27679218474852
$ws.Range("G2").Value = 3.733460641871543
$ws.Range("I2").Value = 33.13978354058633
$ws.Range("J2").Value = 9.635956075275796
$ws.Range("L2").Value = 11.73619395249847
$ws.Range("M2").Value = 19.37031372650769
$ws.Range("N2").Value = 21.20273523982256
$ws.Range("B3").Value = 22.36139988102607
$ws.Range("C3").Value = 10.60700482395662
$ws.Range("E3").Value = 10.2876800664619
$ws.Range("F3").Value = 46.22764367341238
$ws.Range("G3").Value = 3.737427886411731
$ws.Range("I3").Value = 33.1746926346746
$ws.Range("J3").Value = 9.660666996068327
$ws.Range("L3").Value = 11.74736285473079
$ws.Range("M3").Value = 19.32119438438105
$ws.Range("N3").Value = 21.27294099973652
$ws.Range("B4").Value = 22.16241178435306
$ws.Range("C4").Value = 10.29094830652745
$ws.Range("E4").Value = 10.25942472709585
$ws.Range("F4").Value = 46.20838312660449
$ws.Range("G4").Value = 3.739989981559146
$ws.Range("I4").Value = 33.20342203207819
$ws.Range("J4").Value = 9.676497056680541
$ws.Range("L4").Value = 11.75556946606037
$ws.Range("M4").Value = 19.29481914174692
$ws.Range("N4").Value = 21.31804342365202
$ws.Range("B5").Value = 22.08240710906484
$ws.Range("C5").Value = 10.16016747109926
$ws.Range("E5").Value = 10.24772414268564
$ws.Range("F5").Value = 46.20328091412085
$ws.Range("G5").Value = 3.741065909122216
$ws.Range("I5").Value = 33.21695797436008
$ws.Range("J5").Value = 9.683113933899188
$ws.Range("L5").Value = 11.75925353042314
$ws.Range("M5").Value = 19.28502977258219
$ws.Range("N5").Value = 21.33692609701802
$ws.Range("B6").Value = 22.06919089316958
$ws.Range("C6").Value = 10.13833820988836
$ws.Range("E6").Value = 10.24577000933829
$ws.Range("F6").Value = 46.20259958730426
$ws.Range("G6").Value = 3.741246493475468
$ws.Range("I6").Value = 33.2193158574335
$ws.Range("J6").Value = 9.684222707076447
$ws.Range("L6").Value = 11.7598858072768
$ws.Range("M6").Value = 19.28346235454042
$ws.Range("N6").Value = 21.34009197034123
$ws.Range("B7").Value = 22.16132828633143
$ws.Range("C7").Value = 10.28919228977277
$ws.Range("E7").Value = 10.25926768355788
$ws.Range("F7").Value = 46.20830319519824
$ws.Range("G7").Value = 3.740004362764136
$ws.Range("I7").Value = 33.203597187548
$ws.Range("J7").Value = 9.676585621209281
$ws.Range("L7").Value = 11.75561777398629
$ws.Range("M7").Value = 19.29468322813419
$ws.Range("N7").Value = 21.31829604354091
$ws.Range("B8").Value = 22.57715443026774
$ws.Range("C8").Value = 10.93687814564957
$ws.Range("E8").Value = 10.31719507299434
$ws.Range("F8").Value = 46.25758148267123
$ws.Range("G8").Value = 3.734802429988032
$ws.Range("I8").Value = 33.1503027182352
$ws.Range("J8").Value = 9.644340355272549
$ws.Range("L8").Value = 11.73976541641629
$ws.Range("M8").Value = 19.35259652223548
$ws.Range("N8").Value = 21.22652868743057
$ws.Range("B9").Value = 23.41891177700256
$ws.Range("C9").Value = 12.12975601855404
$ws.Range("E9").Value = 10.4249349088674
$ws.Range("F9").Value = 46.44068649803159
$ws.Range("G9").Value = 3.725597194186181
$ws.Range("I9").Value = 33.10391686239163
$ws.Range("J9").Value = 9.586292756869092
$ws.Range("L9").Value = 11.71935371706452
$ws.Range("M9").Value = 19.49583857479395
$ws.Range("N9").Value = 21.0623525225732
$ws.Range("B10").Value = 24.04717749372744
$ws.Range("C10").Value = 12.95014541542297
$ws.Range("E10").Value = 10.50054638302256
$ws.Range("F10").Value = 46.62760444717896
$ws.Range("G10").Value = 3.719433397302538
$ws.Range("I10").Value = 33.10556613518016
$ws.Range("J10").Value = 9.546761709003752
$ws.Range("L10").Value = 11.7108265277745
$ws.Range("M10").Value = 19.61864548536609
$ws.Range("N10").Value = 20.95127093701731
$ws.Range("B11").Value = 24.33382211724222
$ws.Range("C11").Value = 13.30973848365806
$ws.Range("E11").Value = 10.53417964141804
$ws.Range("F11").Value = 46.7239013336596
$ws.Range("G11").Value = 3.716757803011962
$ws.Range("I11").Value = 33.11412358573945
$ws.Range("J11").Value = 9.529445154925391
$ws.Range("L11").Value = 11.70834312562328
$ws.Range("M11").Value = 19.67819961867776
$ws.Range("N11").Value = 20.90279115907137
$ws.Range("B12").Value = 24.44238218660962
$ws.Range("C12").Value = 13.44385307855335
$ws.Range("E12").Value = 10.54680628847935
$ws.Range("F12").Value = 46.76197338624637
$ws.Range("G12").Value = 3.715762954203736
$ws.Range("I12").Value = 33.11848939074507
$ws.Range("J12").Value = 9.522982916642428
$ws.Range("L12").Value = 11.70760258101277
$ws.Range("M12").Value = 19.70126950978684
$ws.Range("N12").Value = 20.88472695859702
$ws.Range("B13").Value = 24.41900288768333
$ws.Range("C13").Value = 13.4150619294855
$ws.Range("E13").Value = 10.54409178353195
$ws.Range("F13").Value = 46.75370267853432
$ws.Range("G13").Value = 3.715976398823183
$ws.Range("I13").Value = 33.11749905567941
$ws.Range("J13").Value = 9.5243704531327
$ws.Range("L13").Value = 11.70775319512827
$ws.Range("M13").Value = 19.69627816131552
$ws.Range("N13").Value = 20.88860434567573
$ws.Range("B14").Value = 24.34275374323489
$ws.Range("C14").Value = 13.3208138784115
$ws.Range("E14").Value = 10.53522063773148
$ws.Range("F14").Value = 46.72700143325964
$ws.Range("G14").Value = 3.716675589288141
$ws.Range("I14").Value = 33.11446019852435
$ws.Range("J14").Value = 9.528911599169211
$ws.Range("L14").Value = 11.70827819986646
$ws.Range("M14").Value = 19.6800872839111
$ws.Range("N14").Value = 20.90129911913207
$ws.Range("B15").Value = 24.29604762720489
$ws.Range("C15").Value = 13.2628138714214
$ws.Range("E15").Value = 10.52977254406679
$ws.Range("F15").Value = 46.71085492302217
$ws.Range("G15").Value = 3.717106248711519
$ws.Range("I15").Value = 33.11274542053535
$ws.Range("J15").Value = 9.531705556231444
$ws.Range("L15").Value = 11.70862578351554
$ws.Range("M15").Value = 19.67023698700138
$ws.Range("N15").Value = 20.90911329751874
$ws.Range("B16").Value = 24.02845232434971
$ws.Range("C16").Value = 12.92636254896637
$ws.Range("E16").Value = 10.49833303422999
$ws.Range("F16").Value = 46.62153681986544
$ws.Range("G16").Value = 3.719610826477834
$ws.Range("I16").Value = 33.10516422781664
$ws.Range("J16").Value = 9.547906738887518
$ws.Range("L16").Value = 11.71101683841221
$ws.Range("M16").Value = 19.6148266655179
$ws.Range("N16").Value = 20.95448040534768
$ws.Range("B17").Value = 23.86443002368452
$ws.Range("C17").Value = 12.71640015451861
$ws.Range("E17").Value = 10.47885065437423
$ws.Range("F17").Value = 46.5696197523377
$ws.Range("G17").Value = 3.721180094386317
$ws.Range("I17").Value = 33.10251512813671
$ws.Range("J17").Value = 9.558015836085515
$ws.Range("L17").Value = 11.71284055999291
$ws.Range("M17").Value = 19.5817704009589
$ws.Range("N17").Value = 20.98283644726387
$ws.Range("B18").Value = 23.7701743741001
$ws.Range("C18").Value = 12.59435948557309
$ws.Range("E18").Value = 10.46757309967386
$ws.Range("F18").Value = 46.54081988384105
$ws.Range("G18").Value = 3.72209478337115
$ws.Range("I18").Value = 33.10172617766026
$ws.Range("J18").Value = 9.563893077232628
$ws.Range("L18").Value = 11.71402093486677
$ws.Range("M18").Value = 19.56310529537613
$ws.Range("N18").Value = 20.99933927589736
$ws.Range("B19").Value = 23.73827905673154
$ws.Range("C19").Value = 12.55282276121623
$ws.Range("E19").Value = 10.4637423762538
$ws.Range("F19").Value = 46.5312514305962
$ws.Range("G19").Value = 3.722406560927168
$ws.Range("I19").Value = 33.10158515259094
$ws.Range("J19").Value = 9.565893808310383
$ws.Range("L19").Value = 11.71444318524882
$ws.Range("M19").Value = 19.55684574318029
$ws.Range("N19").Value = 21.0049600604872
$ws.Range("B20").Value = 23.88188240103444
$ws.Range("C20").Value = 12.73888382286804
$ws.Range("E20").Value = 10.48093201742699
$ws.Range("F20").Value = 46.57503664845567
$ws.Range("G20").Value = 3.721011792958437
$ws.Range("I20").Value = 33.10272106494932
$ws.Range("J20").Value = 9.556933215165083
$ws.Range("L20").Value = 11.71263282626466
$ws.Range("M20").Value = 19.58525336920712
$ws.Range("N20").Value = 20.9797979118183
$ws.Range("B21").Value = 24.36515042607937
$ws.Range("C21").Value = 13.34855331214983
$ws.Range("E21").Value = 10.53782927995752
$ws.Range("F21").Value = 46.73480075537117
$ws.Range("G21").Value = 3.716469723257688
$ws.Range("I21").Value = 33.11532222758976
$ws.Range("J21").Value = 9.527575176927328
$ws.Range("L21").Value = 11.70811857577114
$ws.Range("M21").Value = 19.6848289801627
$ws.Range("N21").Value = 20.89756238420811
$ws.Range("B22").Value = 24.68102796243648
$ws.Range("C22").Value = 13.73499131214584
$ws.Range("E22").Value = 10.57437652595184
$ws.Range("F22").Value = 46.84857199209523
$ws.Range("G22").Value = 3.713608068344191
$ws.Range("I22").Value = 33.13011750020535
$ws.Range("J22").Value = 9.508942394947493
$ws.Range("L22").Value = 11.70633280718915
$ws.Range("M22").Value = 19.75292064040354
$ws.Range("N22").Value = 20.84553014732635
$ws.Range("B23").Value = 24.51246970280025
$ws.Range("C23").Value = 13.52987029765951
$ws.Range("E23").Value = 10.55492893871679
$ws.Range("F23").Value = 46.78699911155002
$ws.Range("G23").Value = 3.715125648696841
$ws.Range("I23").Value = 33.121620107833
$ws.Range("J23").Value = 9.518836548929322
$ws.Range("L23").Value = 11.70717963393883
$ws.Range("M23").Value = 19.71630740754513
$ws.Range("N23").Value = 20.87314427881854
$ws.Range("B24").Value = 23.87399203926301
$ws.Range("C24").Value = 12.72872309292122
$ws.Range("E24").Value = 10.47999127265007
$ws.Range("F24").Value = 46.57258440459422
$ws.Range("G24").Value = 3.721087843103332
$ws.Range("I24").Value = 33.10262567471515
$ws.Range("J24").Value = 9.557422464369845
$ws.Range("L24").Value = 11.71272633180347
$ws.Range("M24").Value = 19.58367766123458
$ws.Range("N24").Value = 20.98117100890988
$ws.Range("B25").Value = 23.18902021149501
$ws.Range("C25").Value = 11.81632833875873
$ws.Range("E25").Value = 10.39641271365947
$ws.Range("F25").Value = 46.38191938514733
$ws.Range("G25").Value = 3.727981655479989
$ws.Range("I25").Value = 33.11020964315897
$ws.Range("J25").Value = 9.601445663726411
$ws.Range("L25").Value = 11.72373660893299
$ws.Range("M25").Value = 19.45396165529186
$ws.Range("N25").Value = 21.10508566990527
